$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.166.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.405.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("E4").Value = "  +0.52%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.54%  "
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("E8").Value = "  +2.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.413.16"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.109"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.95%  "
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.347"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000171"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.840.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.999.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.413.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.84%  "
$ws.Range("E22").Value = "  +1.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.76%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "588.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0953"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.517.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.86%  "
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("E34").Value = "  +1.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.996"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "153.44"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.373"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.34%  "
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("E43").Value = "  +1.85%  "
$ws.Range("E44").Value = "  +10.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0288"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.36"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.593"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0512"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.83%  "
